$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(
    "Mon Feb 24 21:58:58 EST 2025",
    "Mon Feb 24 22:00:12 EST 2025",
    "Mon Feb 24 22:01:25 EST 2025",
    "Mon Feb 24 22:02:37 EST 2025",
    "Mon Feb 24 22:03:51 EST 2025",
    "Mon Feb 24 22:05:04 EST 2025",
    "Mon Feb 24 22:06:17 EST 2025",
    "Mon Feb 24 22:07:29 EST 2025",
    "Mon Feb 24 22:07:41 EST 2025",
    "Mon Feb 24 22:07:52 EST 2025",
    "Mon Feb 24 22:08:03 EST 2025",
    "Mon Feb 24 22:08:14 EST 2025",
    "Mon Feb 24 22:08:25 EST 2025",
    "Mon Feb 24 22:09:38 EST 2025",
    "Mon Feb 24 22:10:50 EST 2025",
    "Mon Feb 24 22:12:03 EST 2025",
    "Mon Feb 24 22:13:16 EST 2025",
    "Mon Feb 24 22:13:28 EST 2025",
    "Mon Feb 24 22:13:39 EST 2025",
    "Mon Feb 24 22:13:49 EST 2025",
    "Mon Feb 24 22:14:00 EST 2025",
    "Mon Feb 24 22:14:11 EST 2025",
    "Mon Feb 24 22:14:22 EST 2025",
    "Mon Feb 24 22:14:33 EST 2025",
    "Mon Feb 24 22:14:44 EST 2025",
    "Mon Feb 24 22:14:55 EST 2025",
    "Mon Feb 24 22:15:05 EST 2025",
    "Mon Feb 24 22:16:18 EST 2025"
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newDates[$i]
}
